# 2017-07-30-Feeding.xlsx -- "algae-added" sheet updates
#  - row 6 (8/3/17): F6/G6 were left blank (no Count 4 / Count 5 taken) -> "N/A",
#    and the actual volume fed (L6) was 500 mL instead of the planned 450 mL.
#  - new row 7 (8/4/17) recorded, with a note about buckets 15 & 16 being mis-fed.
#  - selection/scroll position updated to reflect where the user ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("algae-added")

# --- Row 6 corrections -------------------------------------------------
$ws.Range("F6").Value = "N/A"
$ws.Range("G6").Value = "N/A"
$ws.Range("L6").Value = 500

# --- New row 7 -----------------------------------------------------------
# Copy row 6's formatting (styles, wrap text, number formats, row height)
# down onto row 7 before filling in the new data.
$ws.Range("A6:O6").Copy()
$ws.Range("A7:O7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Application.CutCopyMode = $false
$ws.Rows.Item(7).RowHeight = 32

$ws.Range("A7").Value = 42951
$ws.Range("B7").Value = "500 mL Ciso, 250 609, 250 Chagra"
$ws.Range("C7").Value = 243
$ws.Range("D7").Value = 231
$ws.Range("E7").Value = 131
$ws.Range("F7").Value = 268
$ws.Range("G7").Value = 159
$ws.Range("H7").Formula = "=AVERAGE(C7:G7)"
$ws.Range("I7").Formula = "=(H7*9)/0.0009"
$ws.Range("J7").Formula = "=15000*50000"
$ws.Range("K7").Formula = "=J7/I7"
$ws.Range("L7").Value = 400
$ws.Range("M7").Formula = "=L7*I7"
$ws.Range("N7").Formula = "=M7/15000"
$ws.Range("O7").Value = "Extremely dense, most likely undercounting. 15 and 16 accidentally fed 400 mL mix + 100 Chagra + 100 609"

# --- View state ------------------------------------------------------------
$excel.Goto($ws.Range("B1"), $true)
$ws.Range("O8").Select()
